# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# Price cells that look like plain numbers (e.g. "212.65", "1.00") are
# prefixed with a leading apostrophe so Excel stores them as literal text
# (matching the workbook's original inline-string cells) instead of
# silently coercing them to numeric values and dropping significant
# trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.538.92"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "1.646.33"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'212.65"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").Value = "'0.530"
$ws.Range("E6").Value = "  +4.02%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D8").Value = "'23.56"
$ws.Range("E8").Value = "  -2.62%  "
$ws.Range("D9").Value = "'0.257"
$ws.Range("E9").Value = "  -1.97%  "
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("D12").Value = "1.879.63"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("D13").Value = "1.645.25"
$ws.Range("E13").Value = "  -1.53%  "
$ws.Range("E14").Value = "  +3.79%  "
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("D16").Value = "'64.52"
$ws.Range("E16").Value = "  -2.59%  "
$ws.Range("D17").Value = "27.502.04"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "'232.08"
$ws.Range("E18").Value = "  -3.85%  "
$ws.Range("D19").Value = "0.0₃0725"
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  -3.72%  "
$ws.Range("D23").Value = "'9.79"
$ws.Range("E23").Value = "  +4.39%  "
$ws.Range("D24").Value = "'2.02"
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("D25").Value = "'148.54"
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("E26").Value = "  -2.88%  "
$ws.Range("E27").Value = "  +1.69%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").Value = "'15.61"
$ws.Range("E29").Value = "  -4.58%  "
$ws.Range("E30").Value = "  -2.89%  "
$ws.Range("E32").Value = "  -0.87%  "
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("D34").Value = "1.424.66"
$ws.Range("E34").Value = "  -2.25%  "
$ws.Range("D35").Value = "'1.60"
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("D38").Value = "'0.889"
$ws.Range("E38").Value = "  -4.40%  "
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").Value = "'0.821"
$ws.Range("E42").Value = "  +3.41%  "
$ws.Range("D43").Value = "'5.55"
$ws.Range("E43").Value = "  +2.26%  "
$ws.Range("E44").Value = "  -1.68%  "
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("D46").Value = "'65.13"
$ws.Range("E46").Value = "  -6.90%  "
$ws.Range("D47").Value = "1.789.22"
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("D48").Value = "'1.68"
$ws.Range("E48").Value = "  -2.41%  "
$ws.Range("D49").Value = "'88.18"
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").Value = "'7.78"
$ws.Range("E51").Value = "  -1.94%  "
